$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A104").Value = 13
$ws.Range("B104").Value = 63.05
$ws.Range("C104").Value = 12
$ws.Range("D104").Value = 300
$ws.Range("E104").Value = 4000
$ws.Range("F104").Value = 'QWTMĆU ŚNB|PŁ AĄCDEĘ FGHIJK LŃOÓRS VXYZŹŻ'
$ws.Range("G104").Value = -1403.88365943578
$ws.Range("H104").Value = 'QWTMĆU ŚNB|PŁ AĄCDEĘ FGHIJK LŃOÓRS VXYZŹŻ'
$ws.Range("I104").Value = -1403.8837

$ws.Range("A105").Value = 31
$ws.Range("B105").Value = 136.66
$ws.Range("C105").Value = 12
$ws.Range("D105").Value = 300
$ws.Range("E105").Value = 4000
$ws.Range("F105").Value = 'ĄXTJGH EZYKŃB ACĆDĘF ILŁMNO ÓPQRSŚ UVWŹŻ|'
$ws.Range("G105").Value = -1403.88365943578
$ws.Range("H105").Value = 'ĄXTJGH EZYKŃB ACĆDĘF ILŁMNO ÓPQRSŚ UVWŹŻ|'
$ws.Range("I105").Value = -1403.8837

$ws.Range("A106").Value = 14
$ws.Range("B106").Value = 71.45999999999999
$ws.Range("C106").Value = 12
$ws.Range("D106").Value = 300
$ws.Range("E106").Value = 4000
$ws.Range("F106").Value = 'ĄĘGÓŃQ CVŁXŚO ABĆDEF HIJKLM NPRSTU WYZŹŻ|'
$ws.Range("G106").Value = -1403.88365943578
$ws.Range("H106").Value = 'WYZŹŻ| ĄĘGÓŃQ ABĆDEF HIJKLM NPRSTU CVŁXŚO'
$ws.Range("I106").Value = -1471.5758
$ws.Range("J106").Value = '30 Jun 2021 16:50:07'

$ws.Range("A107").Value = 72
$ws.Range("B107").Value = 319.64
$ws.Range("C107").Value = 12
$ws.Range("D107").Value = 300
$ws.Range("E107").Value = 4000
$ws.Range("F107").Value = 'PVEĄĆQ DJSNUR ABCĘFG HIKLŁM ŃOÓŚTW XYZŹŻ|'
$ws.Range("G107").Value = -1403.88365943578
$ws.Range("H107").Value = 'PVEĄĆQ DJSNUR ABCĘFG HIKLŁM ŃOÓŚTW XYZŹŻ|'
$ws.Range("I107").Value = -1403.8837
$ws.Range("J107").Value = '30 Jun 2021 16:55:30'

$ws.Range("A108").Value = 16
$ws.Range("B108").Value = 109.32
$ws.Range("C108").Value = 12
$ws.Range("D108").Value = 300
$ws.Range("E108").Value = 6000
$ws.Range("F108").Value = 'DŻXSFU QÓLĘZK AĄBCĆE GHIJŁM NŃOPRŚ TVWYŹ|'
$ws.Range("G108").Value = -1403.88365943578
$ws.Range("H108").Value = 'DŻXSFU QÓLĘZK AĄBCĆE GHIJŁM NŃOPRŚ TVWYŹ|'
$ws.Range("I108").Value = -1403.8837
$ws.Range("J108").Value = '30 Jun 2021 17:13:24'

$ws.Range("A109").Value = 12
$ws.Range("B109").Value = 74.62
$ws.Range("C109").Value = 12
$ws.Range("D109").Value = 300
$ws.Range("E109").Value = 6000
$ws.Range("F109").Value = 'ĄĘFLÓN ŹJBIŚŁ ACĆDEG HKMŃOP QRSTUV WXYZŻ|'
$ws.Range("G109").Value = -1403.88365943578
$ws.Range("H109").Value = 'ĄĘFLÓN ŹJBIŚŁ ACĆDEG HKMŃOP QRSTUV WXYZŻ|'
$ws.Range("I109").Value = -1403.8837
$ws.Range("J109").Value = '30 Jun 2021 17:14:44'

$ws.Range("A110").Value = 19
$ws.Range("B110").Value = 119.66
$ws.Range("C110").Value = 12
$ws.Range("D110").Value = 300
$ws.Range("E110").Value = 6000
$ws.Range("F110").Value = 'FYĄMJS ŹOQŃWĆ ABCDEĘ GHIKLŁ NÓPRŚT UVXZŻ|'
$ws.Range("G110").Value = -1403.88365943578
$ws.Range("H110").Value = 'FYĄMJS ŹOQŃWĆ ABCDEĘ GHIKLŁ NÓPRŚT UVXZŻ|'
$ws.Range("I110").Value = -1403.8837
$ws.Range("J110").Value = '30 Jun 2021 17:16:48'

$ws.Range("A111").Value = 9
$ws.Range("B111").Value = 57.42
$ws.Range("C111").Value = 12
$ws.Range("D111").Value = 300
$ws.Range("E111").Value = 6000
$ws.Range("F111").Value = 'AWŚEVK CÓŁŃGH ĄBĆDĘF IJLMNO PQRSTU XYZŹŻ|'
$ws.Range("G111").Value = -1403.88365943578
$ws.Range("H111").Value = 'AWŚEVK CÓŁŃGH ĄBĆDĘF IJLMNO PQRSTU XYZŹŻ|'
$ws.Range("I111").Value = -1403.8837
$ws.Range("J111").Value = '30 Jun 2021 17:17:50'

$ws.Range("A112").Value = 26
$ws.Range("B112").Value = 186.73
$ws.Range("C112").Value = 12
$ws.Range("D112").Value = 300
$ws.Range("E112").Value = 6000
$ws.Range("F112").Value = 'ÓCVKBM ŃPSYWĄ AĆDEĘF GHIJLŁ NOQRŚT UXZŹŻ|'
$ws.Range("G112").Value = -1403.88365943578
$ws.Range("H112").Value = 'ÓCVKBM ŃPSYWĄ AĆDEĘF GHIJLŁ NOQRŚT UXZŹŻ|'
$ws.Range("I112").Value = -1403.8837
$ws.Range("J112").Value = '30 Jun 2021 17:21:01'

$ws.Range("A113").Value = 36
$ws.Range("B113").Value = 232.77
$ws.Range("C113").Value = 12
$ws.Range("D113").Value = 300
$ws.Range("E113").Value = 6000
$ws.Range("F113").Value = 'HNQŁXS JYKLIĘ AĄBCĆD EFGMŃO ÓPRŚTU VWZŹŻ|'
$ws.Range("G113").Value = -1403.88365943578
$ws.Range("H113").Value = 'HNQŁXS JYKLIĘ AĄBCĆD EFGMŃO ÓPRŚTU VWZŹŻ|'
$ws.Range("I113").Value = -1403.8837
$ws.Range("J113").Value = '30 Jun 2021 17:25:00'

$ws.Range("A114").Value = 10
$ws.Range("B114").Value = 66.14
$ws.Range("C114").Value = 12
$ws.Range("D114").Value = 300
$ws.Range("E114").Value = 6000
$ws.Range("F114").Value = 'EXQĘHA WŃTCRL ĄBĆDFG IJKŁMN OÓPSŚU VYZŹŻ|'
$ws.Range("G114").Value = -1403.88365943578
$ws.Range("H114").Value = 'EXQĘHA WŃTCRL ĄBĆDFG IJKŁMN OÓPSŚU VYZŹŻ|'
$ws.Range("I114").Value = -1403.8837
$ws.Range("J114").Value = '30 Jun 2021 17:26:12'

$ws.Range("A115").Value = 21
$ws.Range("B115").Value = 133.34
$ws.Range("C115").Value = 12
$ws.Range("D115").Value = 300
$ws.Range("E115").Value = 6000
$ws.Range("F115").Value = 'WFYDHK ŹĘBOPS AĄCĆEG IJLŁMN ŃÓQRŚT UVXZŻ|'
$ws.Range("G115").Value = -1403.88365943578
$ws.Range("H115").Value = 'AĄCĆEG IJLŁMN ŃÓQRŚT UVXZŻ| WFYDHK ŹĘBOPS'
$ws.Range("I115").Value = -1403.8837
$ws.Range("J115").Value = '30 Jun 2021 17:28:30'

$ws.Range("A116").Value = 11
$ws.Range("B116").Value = 180.9
$ws.Range("C116").Value = 12
$ws.Range("D116").Value = 300
$ws.Range("E116").Value = 6000
$ws.Range("F116").Value = 'ŚŹKPCI RULŃSX AĄBĆDE ĘFGHJŁ MNOÓQT VWYZŻ|'
$ws.Range("G116").Value = -1403.88365943578
$ws.Range("H116").Value = 'ŚŹKPCI RULŃSX AĄBĆDE ĘFGHJŁ MNOÓQT VWYZŻ|'
$ws.Range("I116").Value = -1403.8837
$ws.Range("J116").Value = '30 Jun 2021 18:23:19'

$ws.Range("A117").Value = 14
$ws.Range("B117").Value = 290.88
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 300
$ws.Range("E117").Value = 6000
$ws.Range("F117").Value = 'SŃŻHEĘ I|ĄĆKL ABCDFG JŁMNOÓ PQRŚTU VWXYZŹ'
$ws.Range("G117").Value = -1403.88365943578
$ws.Range("H117").Value = 'SŃŻHEĘ I|ĄĆKL ABCDFG JŁMNOÓ PQRŚTU VWXYZŹ'
$ws.Range("I117").Value = -1403.8837
$ws.Range("J117").Value = '30 Jun 2021 18:30:05'

$ws.Range("A118").Value = 51
$ws.Range("B118").Value = 330.99
$ws.Range("C118").Value = 12
$ws.Range("D118").Value = 300
$ws.Range("E118").Value = 6000
$ws.Range("F118").Value = 'MŹZLIF ŁAEÓCŻ ĄBĆDĘG HJKNŃO PQRSŚT UVWXY|'
$ws.Range("G118").Value = -1403.88365943578
$ws.Range("H118").Value = 'MŹZLIF ŁAEÓCŻ ĄBĆDĘG HJKNŃO PQRSŚT UVWXY|'
$ws.Range("I118").Value = -1403.8837
$ws.Range("J118").Value = '30 Jun 2021 18:48:10'

$ws.Range("A119").Value = 12
$ws.Range("B119").Value = 84.90000000000001
$ws.Range("C119").Value = 12
$ws.Range("D119").Value = 300
$ws.Range("E119").Value = 6000
$ws.Range("F119").Value = 'SPMORŃ IX|ĆKB AĄCDEĘ FGHJLŁ NÓQŚTU VWYZŹŻ'
$ws.Range("G119").Value = -1403.88365943578
$ws.Range("H119").Value = 'SPMORŃ IX|ĆKB AĄCDEĘ FGHJLŁ NÓQŚTU VWYZŹŻ'
$ws.Range("I119").Value = -1403.8837
$ws.Range("J119").Value = '30 Jun 2021 18:54:16'

$ws.Range("A120").Value = 15
$ws.Range("B120").Value = 108.26
$ws.Range("C120").Value = 12
$ws.Range("D120").Value = 300
$ws.Range("E120").Value = 6000
$ws.Range("F120").Value = 'KĘFPWC ÓNALXV ĄBĆDEG HIJŁMŃ OQRSŚT UYZŹŻ|'
$ws.Range("G120").Value = -1403.88365943578
$ws.Range("H120").Value = 'KĘFPWC ÓNALXV ĄBĆDEG HIJŁMŃ OQRSŚT UYZŹŻ|'
$ws.Range("I120").Value = -1403.8837
$ws.Range("J120").Value = '30 Jun 2021 18:56:09'

$ws.Range("A121").Value = 51
$ws.Range("B121").Value = 357.68
$ws.Range("C121").Value = 12
$ws.Range("D121").Value = 300
$ws.Range("E121").Value = 6000
$ws.Range("F121").Value = 'IŁNÓUŹ JMCĄAĘ BĆDEFG HKLŃOP QRSŚTV WXYZŻ|'
$ws.Range("G121").Value = -1403.88365943578
$ws.Range("H121").Value = 'IŁNÓUŹ JMCĄAĘ QRSŚTV BĆDEFG HKLŃOP WXYZŻ|'
$ws.Range("I121").Value = -1453.0834
$ws.Range("J121").Value = '30 Jun 2021 19:02:11'

$ws.Range("A122").Value = 15
$ws.Range("B122").Value = 95.59
$ws.Range("C122").Value = 12
$ws.Range("D122").Value = 300
$ws.Range("E122").Value = 6000
$ws.Range("F122").Value = 'CŹĄEIG ŁHRPBA ĆDĘFJK LMNŃOÓ QSŚTUV WXYZŻ|'
$ws.Range("G122").Value = -1403.88365943578
$ws.Range("H122").Value = 'CŹĄEIG ŁHRPBA ĆDĘFJK LMNŃOÓ QSŚTUV WXYZŻ|'
$ws.Range("I122").Value = -1403.8837
$ws.Range("J122").Value = '30 Jun 2021 19:04:36'

$ws.Range("A123").Value = 18
$ws.Range("B123").Value = 111.65
$ws.Range("C123").Value = 12
$ws.Range("D123").Value = 300
$ws.Range("E123").Value = 6000
$ws.Range("F123").Value = 'LACXGJ HVÓO|Ć ĄBDEĘF IKŁMNŃ PQRSŚT UWYZŹŻ'
$ws.Range("G123").Value = -1403.88365943578
$ws.Range("H123").Value = 'LACXGJ HVÓO|Ć ĄBDEĘF IKŁMNŃ PQRSŚT UWYZŹŻ'
$ws.Range("I123").Value = -1403.8837
$ws.Range("J123").Value = '30 Jun 2021 19:06:32'
